$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (column-header style row, mirrors the ACCESSION table headers) ---
$ws.Range("A2").Value = "ACCESSION NO"
$ws.Range("C2").Value = "TITLE"
$ws.Range("D2").Value = "ACQUISITION NO."
$ws.Range("E2").Value = "Item"
$ws.Range("F2").Value = "ITEM DESCRIPTION"
$ws.Range("G2").Value = "LOCATION | SECTION"
$ws.Range("K2").Value = "QTY"
$ws.Range("L2").Value = "AR NUMBER"
$ws.Range("M2").Value = "2022/23 RFID Number"
$ws.Range("N2").Value = "COLLECTIONS"

# --- Row 3 ---
$ws.Range("A3").Value = "EW203-1-1"
$ws.Range("C3").Value = "AFRICAN SCHOOL FEEDING SCHEME"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "01/04/2014"
$ws.Range("E3").Value = "Item"
$ws.Range("F3").Value = "NEGATIVE B/W MEDIUM FORMAT"
$ws.Range("G3").Value = "S-3D (EW 203-208/ EW 210,217/ EW 220,221 & 223) | 3.25"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = "421V-PH"
$ws.Range("N3").Value = "HAIN, Peter"

# --- Row 4 ---
$ws.Range("A4").Value = "EW204-1-1"
$ws.Range("C4").Value = "ANC SCHOOL, GERMEISTON, ANTI BANTU EDUCATION"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "01/04/2014"
$ws.Range("E4").Value = "Item"
$ws.Range("F4").Value = "NEGATIVE B/W MEDIUM FORMAT"
$ws.Range("G4").Value = "S-3D (EW 203-208/ EW 210,217/ EW 220,221 & 223) | 3.25"
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = "869V-PH"
$ws.Range("N4").Value = "HAIN, Peter"

# --- Row 5 ---
$ws.Range("A5").Value = "EW205-1-1"
$ws.Range("C5").Value = "KOAKOA PASS/ APARTHEID PHONE BOOTHS/ SOWETO"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "01/04/2014"
$ws.Range("E5").Value = "Item"
$ws.Range("F5").Value = "NEGATIVE B/W MEDIUM FORMAT"
$ws.Range("G5").Value = "S-3D (EW 203-208/ EW 210,217/ EW 220,221 & 223) | 3.25"
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = "14370V-PH"
$ws.Range("N5").Value = "HAIN, Peter"

# --- Row 6 ---
$ws.Range("A6").Value = "EW206-1-3"
$ws.Range("C6").Value = "HOUSING - CAPE FISHERMEN CA 1936"
$ws.Range("D6").Value = "date in title"
$ws.Range("E6").Value = "Item"
$ws.Range("F6").Value = "NEGATIVE B/W MEDIUM FORMAT"
$ws.Range("G6").Value = "S-3D (EW 203-208/ EW 210,217/ EW 220,221 & 223) | 3.25"
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = "4797V-PH"
$ws.Range("N6").Value = "HAIN, Peter"

# --- Hide the helper columns K:N that hold the QTY/AR NUMBER/RFID/COLLECTIONS data ---
$ws.Range("K1:N1").EntireColumn.Hidden = $true
$ws.Range("K1:N1").EntireColumn.ColumnWidth = 0
